$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Statistics (L14)")
$ws.Range("B3").Value = 245.5714285714286
$ws.Range("D3").Value = 19.92857142857143
$ws.Range("E3").Value = 0.2142857142857143
$ws.Range("F3").Value = 10602.92857142857
$ws.Range("G3").Value = 0.04356462792655074
$ws.Range("H3").Value = 922.7142857142857
$ws.Range("J3").Value = 3.784982313396027
$ws.Range("K3").Value = 669.7142857142857
$ws.Range("L3").Value = 1.686480578325963
$ws.Range("M3").Value = 28.35714285714286
$ws.Range("N3").Value = -0.1428571428571428
$ws.Range("B4").Value = 62.06288905458879
$ws.Range("D4").Value = 17.28255806259054
$ws.Range("E4").Value = 0.4258153136263202
$ws.Range("F4").Value = 1030.92829447622
$ws.Range("G4").Value = 0.3197112084771773
$ws.Range("J4").Value = 2.764455872658345
$ws.Range("L4").Value = 5.831209625978868
$ws.Range("M4").Value = 0.9287827316640653
$ws.Range("N4").Value = 0.3631365196012815
$ws.Range("F5").Value = 9083
$ws.Range("H5").Value = 860
$ws.Range("K5").Value = 533
$ws.Range("F6").Value = 9782
$ws.Range("G6").Value = -0.16375
$ws.Range("H6").Value = 873.5
$ws.Range("K6").Value = 653.5
$ws.Range("L6").Value = -0.5022611644997174
$ws.Range("M6").Value = 28
$ws.Range("N6").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 10513.5
$ws.Range("H7").Value = 923
$ws.Range("K7").Value = 671
$ws.Range("L7").Value = -0.1920289855072463
$ws.Range("M7").Value = 29
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 11482.25
$ws.Range("H8").Value = 972.25
$ws.Range("J8").Value = 5.634268324623655
$ws.Range("K8").Value = 707.5
$ws.Range("M8").Value = 29
$ws.Range("N8").Value = 0
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 12197
$ws.Range("H9").Value = 989
$ws.Range("K9").Value = 769
$ws.Range("M9").Value = 29
$ws.Range("N9").Value = 0

$ws = $wb.Worksheets.Item("Statistics (All)")
$ws.Range("B3").Value = 145.202380952381
$ws.Range("C3").Value = 11.77380952380952
$ws.Range("D3").Value = 9.154761904761905
$ws.Range("E3").Value = 0.3452380952380952
$ws.Range("F3").Value = 3871.595238095238
$ws.Range("G3").Value = 0.3662432723157148
$ws.Range("H3").Value = 374.9761904761905
$ws.Range("I3").Value = 0.6605708486938211
$ws.Range("J3").Value = 7.317510264686279
$ws.Range("K3").Value = 169.3690476190476
$ws.Range("L3").Value = 0.9412141240030396
$ws.Range("M3").Value = 10.26190476190476
$ws.Range("N3").Value = -0.1349206349206349
$ws.Range("B4").Value = 109.9938420797618
$ws.Range("C4").Value = 14.38234459033328
$ws.Range("D4").Value = 15.28326437061617
$ws.Range("E4").Value = 0.7196727438845412
$ws.Range("F4").Value = 3930.157316570051
$ws.Range("G4").Value = 1.350534704742466
$ws.Range("H4").Value = 371.3236816238768
$ws.Range("I4").Value = 2.135947705410573
$ws.Range("J4").Value = 7.65604245748536
$ws.Range("K4").Value = 248.469473035235
$ws.Range("L4").Value = 4.641660121019395
$ws.Range("M4").Value = 11.34573753369923
$ws.Range("N4").Value = 0.3617192607698794
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 309.25
$ws.Range("G6").Value = -0.2206060606060606
$ws.Range("L6").Value = -0.175
$ws.Range("F7").Value = 2460
$ws.Range("J7").Value = 4.714285714285714
$ws.Range("B8").Value = 231.5
$ws.Range("C8").Value = 18.25
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 7011.25
$ws.Range("I8").Value = 0.78125
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 12197
$ws.Range("H9").Value = 989
$ws.Range("K9").Value = 769
$ws.Range("M9").Value = 29
$ws.Range("N9").Value = 1

$ws = $wb.Worksheets.Item("Kosovo Raw Data")
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 1.5
$ws.Range("C4").Value = 10
$ws.Range("C5").Value = 16
$ws.Range("C6").Value = 23
$ws.Range("C7").Value = 31
$ws.Range("C8").Value = 33
$ws.Range("C9").Value = 37
$ws.Range("C10").Value = 44
$ws.Range("C11").Value = 52
$ws.Range("C12").Value = 53
$ws.Range("C13").Value = 57
$ws.Range("C14").Value = 62
$ws.Range("C15").Value = 76
$ws.Range("C16").Value = 83
$ws.Range("C17").Value = 94
$ws.Range("C18").Value = 120
$ws.Range("C19").Value = 219
$ws.Range("C20").Value = 246
$ws.Range("C21").Value = 270
$ws.Range("C22").Value = 283
$ws.Range("C23").Value = 318
$ws.Range("C24").Value = 339
$ws.Range("C25").Value = 416
$ws.Range("C26").Value = 566
$ws.Range("C27").Value = 576
$ws.Range("C28").Value = 689
$ws.Range("C29").Value = 731
$ws.Range("C30").Value = 765
$ws.Range("C31").Value = 884
$ws.Range("C32").Value = 990
$ws.Range("C33").Value = 1036
$ws.Range("C34").Value = 1081
$ws.Range("C35").Value = 1150
$ws.Range("C36").Value = 1266
$ws.Range("C37").Value = 1379
$ws.Range("C38").Value = 1499
$ws.Range("C39").Value = 1742
$ws.Range("C40").Value = 1870
$ws.Range("C41").Value = 1963
$ws.Range("C42").Value = 2092
$ws.Range("C43").Value = 2351
$ws.Range("C44").Value = 2569
$ws.Range("C45").Value = 2686
$ws.Range("C46").Value = 2875
$ws.Range("C47").Value = 3005
$ws.Range("C48").Value = 3255
$ws.Range("C49").Value = 3409
$ws.Range("C50").Value = 3546
$ws.Range("B51").Value = 172
$ws.Range("C51").Value = 3718
$ws.Range("D51").Value = 0.2554744525547445
$ws.Range("H51").Value = 20.93023255813954
$ws.Range("C52").Value = 3913
$ws.Range("D52").Value = 0.1337209302325582
$ws.Range("C53").Value = 4122
$ws.Range("C54").Value = 4364
$ws.Range("B55").Value = 294
$ws.Range("C55").Value = 4658
$ws.Range("D55").Value = 0.2148760330578512
$ws.Range("H55").Value = 17.3469387755102
$ws.Range("C56").Value = 5011
$ws.Range("D56").Value = 0.2006802721088434
$ws.Range("C57").Value = 5202
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 15
$ws.Range("N57").Value = -1
$ws.Range("C58").Value = 5374
$ws.Range("L58").Value = 3
$ws.Range("N58").Value = 0
$ws.Range("C59").Value = 5629
$ws.Range("N59").Value = -0.6666666666666667
$ws.Range("C60").Value = 5892
$ws.Range("B61").Value = 291
$ws.Range("C61").Value = 6183
$ws.Range("D61").Value = 0.1064638783269962
$ws.Range("H61").Value = 9.621993127147768
$ws.Range("C62").Value = 6480
$ws.Range("D62").Value = 0.02061855670103085
$ws.Range("B63").Value = 225
$ws.Range("C63").Value = 6705
$ws.Range("D63").Value = -0.2424242424242424
$ws.Range("H63").Value = 7.555555555555555
$ws.Range("C64").Value = 6955
$ws.Range("D64").Value = 0.1111111111111112
$ws.Range("C65").Value = 7180
$ws.Range("B66").Value = 242
$ws.Range("C66").Value = 7422
$ws.Range("D66").Value = 0.0755555555555556
$ws.Range("H66").Value = 2.892561983471075
$ws.Range("C67").Value = 7643
$ws.Range("D67").Value = -0.08677685950413228
$ws.Range("B68").Value = 230
$ws.Range("C68").Value = 7873
$ws.Range("D68").Value = 0.04072398190045257
$ws.Range("H68").Value = 4.347826086956522
$ws.Range("C69").Value = 8285
$ws.Range("D69").Value = 0.7913043478260871
$ws.Range("E69").Value = 28
$ws.Range("F69").Value = 851
$ws.Range("G69").Value = 1.8
$ws.Range("H69").Value = 6.796116504854369
$ws.Range("C70").Value = 8539
$ws.Range("F70").Value = 855
$ws.Range("G70").Value = -0.8571428571428572
$ws.Range("I70").Value = 22
$ws.Range("J70").Value = 403
$ws.Range("K70").Value = -0.5111111111111111
$ws.Range("L70").Value = 1
$ws.Range("M70").Value = 26
$ws.Range("N70").Value = -0.6666666666666667
$ws.Range("C71").Value = 8759
$ws.Range("F71").Value = 856
$ws.Range("I71").Value = 87
$ws.Range("J71").Value = 490
$ws.Range("K71").Value = 2.954545454545455
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = -1
$ws.Range("C72").Value = 9083
$ws.Range("F72").Value = 860
$ws.Range("I72").Value = 43
$ws.Range("J72").Value = 533
$ws.Range("K72").Value = -0.5057471264367817
$ws.Range("N72").Value = 0
$ws.Range("C73").Value = 9312
$ws.Range("F73").Value = 861
$ws.Range("J73").Value = 561
$ws.Range("K73").Value = -0.3488372093023255
$ws.Range("L73").Value = 1
$ws.Range("M73").Value = 27
$ws.Range("C74").Value = 9555
$ws.Range("F74").Value = 862
$ws.Range("J74").Value = 622
$ws.Range("M74").Value = 28
$ws.Range("C75").Value = 9732
$ws.Range("F75").Value = 870
$ws.Range("J75").Value = 653
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 28
$ws.Range("N75").Value = -1
$ws.Range("C76").Value = 9932
$ws.Range("F76").Value = 884
$ws.Range("J76").Value = 655
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 28
$ws.Range("N76").Value = 0
$ws.Range("B77").Value = 197
$ws.Range("C77").Value = 10129
$ws.Range("D77").Value = -0.01500000000000001
$ws.Range("F77").Value = 895
$ws.Range("H77").Value = 5.583756345177665
$ws.Range("J77").Value = 657
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 28
$ws.Range("N77").Value = 0
$ws.Range("C78").Value = 10386
$ws.Range("D78").Value = 0.3045685279187818
$ws.Range("F78").Value = 919
$ws.Range("J78").Value = 671
$ws.Range("L78").Value = 1
$ws.Range("M78").Value = 29
$ws.Range("N78").Value = 0
$ws.Range("C79").Value = 10641
$ws.Range("F79").Value = 927
$ws.Range("J79").Value = 671
$ws.Range("M79").Value = 29
$ws.Range("C80").Value = 10941
$ws.Range("F80").Value = 945
$ws.Range("J80").Value = 690
$ws.Range("M80").Value = 29
$ws.Range("C81").Value = 11177
$ws.Range("F81").Value = 955
$ws.Range("J81").Value = 691
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 29
$ws.Range("C82").Value = 11584
$ws.Range("F82").Value = 978
$ws.Range("J82").Value = 713
$ws.Range("M82").Value = 29
$ws.Range("N82").Value = 0
$ws.Range("C83").Value = 11786
$ws.Range("F83").Value = 985
$ws.Range("J83").Value = 736
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 29
$ws.Range("C84").Value = 11986
$ws.Range("F84").Value = 988
$ws.Range("J84").Value = 754
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 29
$ws.Range("N84").Value = 0
$ws.Range("C85").Value = 12197
$ws.Range("F85").Value = 989
$ws.Range("J85").Value = 769
$ws.Range("M85").Value = 29
$ws.Range("N85").Value = 0
